$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.313.93'
$ws.Range("E2").Value = '  +0.49%  '

$ws.Range("D3").Value = '3.801.35'
$ws.Range("E3").Value = '  +1.38%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '604.26'
$ws.Range("E5").Value = '  +0.25%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '165.36'
$ws.Range("E6").Value = '  -2.39%  '

$ws.Range("D7").Value = '3.802.10'
$ws.Range("E7").Value = '  +1.41%  '

$ws.Range("E8").Value = '  -0.02%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.539'
$ws.Range("E9").Value = '  +0.84%  '

$ws.Range("E10").Value = '  +3.40%  '

$ws.Range("E11").Value = '  -0.14%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.462'
$ws.Range("E12").Value = '  -0.17%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000247'
$ws.Range("E14").Value = '  -0.30%  '

$ws.Range("D15").Value = '4.437.13'
$ws.Range("E15").Value = '  +1.37%  '

$ws.Range("D16").Value = '3.801.76'
$ws.Range("E16").Value = '  +1.07%  '

$ws.Range("D17").Value = '69.425.36'
$ws.Range("E17").Value = '  +0.70%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.47'
$ws.Range("E18").Value = '  +2.46%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.53'
$ws.Range("E19").Value = '  +2.82%  '

$ws.Range("E20").Value = '  -0.34%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.28'
$ws.Range("E21").Value = '  +4.80%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '493.85'
$ws.Range("E22").Value = '  -0.65%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.725'
$ws.Range("E23").Value = '  -0.30%  '

$ws.Range("E24").Value = '  -1.78%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '84.87'
$ws.Range("E25").Value = '  -0.71%  '

$ws.Range("E26").Value = '  -2.43%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.31'
$ws.Range("E27").Value = '  -0.53%  '

$ws.Range("E28").Value = '  -1.94%  '

$ws.Range("E29").Value = '  +0.10%  '

$ws.Range("E30").Value = '  +0.25%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.12'
$ws.Range("E31").Value = '  +2.35%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.41'
$ws.Range("E32").Value = '  -4.92%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '32.12'
$ws.Range("E33").Value = '  +0.70%  '

$ws.Range("D34").Value = '3.946.90'
$ws.Range("E34").Value = '  +1.35%  '

$ws.Range("D35").Value = '3.748.09'
$ws.Range("E35").Value = '  +1.74%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.107'
$ws.Range("E36").Value = '  -1.03%  '

$ws.Range("E37").Value = '  +6.23%  '

$ws.Range("E38").Value = '  +0.45%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.96'
$ws.Range("E39").Value = '  +1.65%  '

$ws.Range("E40").Value = '  -0.01%  '

$ws.Range("E41").Value = '  +0.09%  '

$ws.Range("E42").Value = '  +3.27%  '

$ws.Range("E43").Value = '  +0.83%  '

$ws.Range("B44").Value = 'OKB'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '48.46'
$ws.Range("E44").Value = '  -0.97%  '

$ws.Range("B45").Value = 'Bittensor'
$ws.Range("C45").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '425.77'
$ws.Range("E45").Value = '  -2.87%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.43'
$ws.Range("E46").Value = '  -0.50%  '

$ws.Range("E47").Value = '  +0.00%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '142.36'
$ws.Range("E48").Value = '  +0.78%  '

$ws.Range("D49").Value = '2.818.53'
$ws.Range("E49").Value = '  +1.33%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '39.93'
$ws.Range("E50").Value = '  -1.75%  '
